$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post at row 448 ("「カタールが好き」") was removed from the source list.
# Deleting the entire row shifts every subsequent row up by one, matching
# the diff (dimension shrinks from A1:C586 to A1:C585, and every row below
#448 is renumbered down by one with otherwise identical content).
$ws.Rows.Item(448).Delete()
